$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.446.68'
$ws.Range('E2').Value = '  +0.62%  '
$ws.Range('D3').Value = '1.697.75'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.84'
$ws.Range('E5').Value = '  +0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5494'
$ws.Range('E6').Value = '  +4.39%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2739'
$ws.Range('E8').Value = '  +1.56%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06453'
$ws.Range('E9').Value = '  +0.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.00'
$ws.Range('E10').Value = '  -0.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07691'
$ws.Range('E11').Value = '  +2.35%  '
$ws.Range('D12').Value = '1.716.87'
$ws.Range('E12').Value = '  +2.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.555'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5863'
$ws.Range('E14').Value = '  +1.08%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008420'
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.76'
$ws.Range('E16').Value = '  +2.38%  '
$ws.Range('D17').Value = '26.491.01'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.953'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.00'
$ws.Range('E20').Value = '  +1.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '191.85'
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.264'
$ws.Range('E22').Value = '  +1.02%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.03'
$ws.Range('E24').Value = '  +2.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1320'
$ws.Range('E25').Value = '  +6.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.910'
$ws.Range('E26').Value = '  +2.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.86'
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06261'
$ws.Range('E28').Value = '  -4.55%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.381'
$ws.Range('E29').Value = '  +1.70%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.332'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.615'
$ws.Range('E31').Value = '  +1.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.598'
$ws.Range('E32').Value = '  +0.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.692'
$ws.Range('E33').Value = '  +2.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.041'
$ws.Range('E34').Value = '  +1.56%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6177'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.409'
$ws.Range('E36').Value = '  +0.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.761'
$ws.Range('E37').Value = '  +2.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01649'
$ws.Range('E38').Value = '  +1.79%  '
$ws.Range('D39').Value = '1.118.39'
$ws.Range('E39').Value = '  +1.17%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.131'
$ws.Range('E40').Value = '  -3.88%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8780'
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('E42').Value = '  +0.12%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.30'
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('D44').Value = '1.849.44'
$ws.Range('E44').Value = '  +1.20%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.63'
$ws.Range('E45').Value = '  +1.49%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000109'
$ws.Range('E46').Value = '  -4.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.214'
$ws.Range('E47').Value = '  +0.83%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.009'
$ws.Range('E48').Value = '  +0.32%  '
$ws.Range('E49').Value = '  +0.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.130'
$ws.Range('E50').Value = '  +1.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4300'
$ws.Range('E51').Value = '  -0.10%  '
